# Forms the consolidated report: recompute the "Absent" column (H) from the
# "Real" column (E). A day counts as Absent (H=1) whenever the student was
# not marked as a "Real" attendance (E=0); otherwise Absent is 0.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 3) { $lastRow = 21 }

for ($r = 3; $r -le $lastRow; $r++) {
    $real = $ws.Range("E$r").Value2
    if ($real -eq $null) { $real = 0 }
    if ($real -eq 0) {
        $ws.Range("H$r").Value = 1
    } else {
        $ws.Range("H$r").Value = 0
    }
}
